# Apply updated "dSF" (column F) values to Sheet1, row by row.
# These values were repulled/recalculated from source data (per commit message:
# "repull data, push all data, mean calculation"); only column F changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    4  = -3
    5  = 1
    6  = 3
    7  = 4
    9  = 4
    10 = -2
    11 = 8
    12 = 4
    13 = -2
    14 = -1
    16 = 1
    19 = 2
    20 = 4
    21 = -1
    22 = -1
    23 = 3
    28 = 0
    31 = -1
    34 = -3
    35 = 1
    36 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
